# Swap the xDegrees (col B) and yDegrees (col C) values for each data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 17; $row++) {
    $bCell = $ws.Cells.Item($row, 2)
    $cCell = $ws.Cells.Item($row, 3)

    $bValue = $bCell.Value2
    $cValue = $cCell.Value2

    # Force text storage (these columns hold numeric-looking text, e.g.
    # coordinates, not real numbers) before assigning so Excel does not
    # auto-coerce the swapped strings into numeric cells.
    $bCell.NumberFormat = "@"
    $cCell.NumberFormat = "@"

    $bCell.Value = $cValue
    $cCell.Value = $bValue
}
